# netCrypto.xlsx update: refresh the daily T2 figure and move the
# selection cursor to T2 (matching the new active cell in the sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")
$ws.Activate()

$ws.Range("T2").Value = 592482
$ws.Range("T2").Select()
